$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "10+45=55"
$t.Cell(1, 2).Range.Text = "43-16=27"
$t.Cell(1, 3).Range.Text = "76-59=17"
$t.Cell(1, 4).Range.Text = "57-13=44"
$t.Cell(1, 5).Range.Text = "5-4=1"

$t.Cell(2, 1).Range.Text = "85+4=89"
$t.Cell(2, 2).Range.Text = "24+57=81"
$t.Cell(2, 3).Range.Text = "97-75=22"
$t.Cell(2, 4).Range.Text = "80-27=53"
$t.Cell(2, 5).Range.Text = "25+1=26"

$t.Cell(3, 1).Range.Text = "89-1=88"
$t.Cell(3, 2).Range.Text = "96-25=71"
$t.Cell(3, 3).Range.Text = "61+26=87"
$t.Cell(3, 4).Range.Text = "91-55=36"
$t.Cell(3, 5).Range.Text = "61+14=75"

$t.Cell(4, 1).Range.Text = "10+55=65"
$t.Cell(4, 2).Range.Text = "63+2=65"
$t.Cell(4, 3).Range.Text = "9-0=9"
$t.Cell(4, 4).Range.Text = "23+21=44"
$t.Cell(4, 5).Range.Text = "87-28=59"

$t.Cell(5, 1).Range.Text = "93-5=88"
$t.Cell(5, 2).Range.Text = "1+59=60"
$t.Cell(5, 3).Range.Text = "97-45=52"
$t.Cell(5, 4).Range.Text = "11+41=52"
$t.Cell(5, 5).Range.Text = "54+23=77"

$t.Cell(6, 1).Range.Text = "67+12=79"
$t.Cell(6, 2).Range.Text = "29-14=15"
$t.Cell(6, 3).Range.Text = "24+9=33"
$t.Cell(6, 4).Range.Text = "88-32=56"
$t.Cell(6, 5).Range.Text = "38+25=63"

$t.Cell(7, 1).Range.Text = "55+41=96"
$t.Cell(7, 2).Range.Text = "4-0=4"
$t.Cell(7, 3).Range.Text = "32+57=89"
$t.Cell(7, 4).Range.Text = "38+27=65"
$t.Cell(7, 5).Range.Text = "16-4=12"

$t.Cell(8, 1).Range.Text = "93-63=30"
$t.Cell(8, 2).Range.Text = "68+0=68"
$t.Cell(8, 3).Range.Text = "76-72=4"
$t.Cell(8, 4).Range.Text = "55+37=92"
$t.Cell(8, 5).Range.Text = "40+1=41"

$t.Cell(9, 1).Range.Text = "8+66=74"
$t.Cell(9, 2).Range.Text = "41-36=5"
$t.Cell(9, 3).Range.Text = "77-48=29"
$t.Cell(9, 4).Range.Text = "25+43=68"
$t.Cell(9, 5).Range.Text = "82-6=76"

$t.Cell(10, 1).Range.Text = "82+17=99"
$t.Cell(10, 2).Range.Text = "45-24=21"
$t.Cell(10, 3).Range.Text = "16-2=14"
$t.Cell(10, 4).Range.Text = "92-50=42"
$t.Cell(10, 5).Range.Text = "30-1=29"

$t.Cell(11, 1).Range.Text = "29+28=57"
$t.Cell(11, 2).Range.Text = "6+5=11"
$t.Cell(11, 3).Range.Text = "27-4=23"
$t.Cell(11, 4).Range.Text = "58-4=54"
$t.Cell(11, 5).Range.Text = "40+56=96"

$t.Cell(12, 1).Range.Text = "43-25=18"
$t.Cell(12, 2).Range.Text = "5-0=5"
$t.Cell(12, 3).Range.Text = "71-64=7"
$t.Cell(12, 4).Range.Text = "64-48=16"
$t.Cell(12, 5).Range.Text = "23-7=16"

$t.Cell(13, 1).Range.Text = "26+70=96"
$t.Cell(13, 2).Range.Text = "3+17=20"
$t.Cell(13, 3).Range.Text = "49-19=30"
$t.Cell(13, 4).Range.Text = "55+32=87"
$t.Cell(13, 5).Range.Text = "36+63=99"

$t.Cell(14, 1).Range.Text = "98-73=25"
$t.Cell(14, 2).Range.Text = "61-18=43"
$t.Cell(14, 3).Range.Text = "13+38=51"
$t.Cell(14, 4).Range.Text = "77-25=52"
$t.Cell(14, 5).Range.Text = "83-55=28"

$t.Cell(15, 1).Range.Text = "33+10=43"
$t.Cell(15, 2).Range.Text = "27+6=33"
$t.Cell(15, 3).Range.Text = "77-33=44"
$t.Cell(15, 4).Range.Text = "13+1=14"
$t.Cell(15, 5).Range.Text = "46+29=75"

$t.Cell(16, 1).Range.Text = "46-46=0"
$t.Cell(16, 2).Range.Text = "70+20=90"
$t.Cell(16, 3).Range.Text = "68+10=78"
$t.Cell(16, 4).Range.Text = "71-40=31"
$t.Cell(16, 5).Range.Text = "74+16=90"

$t.Cell(17, 1).Range.Text = "37-26=11"
$t.Cell(17, 2).Range.Text = "51-2=49"
$t.Cell(17, 3).Range.Text = "7+88=95"
$t.Cell(17, 4).Range.Text = "18+25=43"
$t.Cell(17, 5).Range.Text = "62+24=86"

$t.Cell(18, 1).Range.Text = "22+36=58"
$t.Cell(18, 2).Range.Text = "88-12=76"
$t.Cell(18, 3).Range.Text = "39+11=50"
$t.Cell(18, 4).Range.Text = "34+50=84"
$t.Cell(18, 5).Range.Text = "46-45=1"

$t.Cell(19, 1).Range.Text = "76-6=70"
$t.Cell(19, 2).Range.Text = "77-23=54"
$t.Cell(19, 3).Range.Text = "44+3=47"
$t.Cell(19, 4).Range.Text = "7+55=62"
$t.Cell(19, 5).Range.Text = "90-74=16"

$t.Cell(20, 1).Range.Text = "91-57=34"
$t.Cell(20, 2).Range.Text = "9+46=55"
$t.Cell(20, 3).Range.Text = "41-10=31"
$t.Cell(20, 4).Range.Text = "55-17=38"
$t.Cell(20, 5).Range.Text = "87-10=77"

